# Edit for PubType.xlsx
# The "user" column (originally column D) moves to become column B,
# pushing "pub_type_name" (old B) and "approved" (old C) one column to
# the right (becoming C and D respectively).
#
# This is implemented as:
#   1. Insert a blank column before column B, shifting B:D -> C:E
#      (shift-right so column widths/formatting travel with the data)
#   2. Cut the "user" data now sitting in column E and paste it into
#      the newly inserted column B (this preserves original cell types,
#      e.g. keeps text "TRUE" as text instead of Excel auto-converting
#      it to a boolean)
#   3. Delete the now-empty column E, shifting everything back to the left
#   4. Update the active cell / selection to match the authored state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftToRight = -4161
$xlShiftToLeft  = -4159

# 1. Make room: push pub_type_name/approved/user one column to the right
$ws.Range("B1:B14").Insert($xlShiftToRight)

# 2. Move "user" (now in column E) into the freshly inserted column B
$ws.Range("E1:E14").Cut($ws.Range("B1:B14"))

# 3. Remove the now-empty column E, shifting remaining cells left
$ws.Range("E1:E14").Delete($xlShiftToLeft)

# 4. Restore the selection / active cell as recorded in the workbook
$ws.Range("F5").Select()
